# Apply the changes described by the commit:
# "regenerate instance to have positive average demands during the last periods"

$wb = $excel.ActiveWorkbook

# --- Productdata sheet: update Leadtimes (column G) for rows 2-5 ---
$wsProductdata = $wb.Worksheets.Item("Productdata")
$wsProductdata.Range("G2").Value = 49
$wsProductdata.Range("G3").Value = 21
$wsProductdata.Range("G4").Value = 35
$wsProductdata.Range("G5").Value = 70

# Keep the blank placeholder column (H) genuinely blank - re-assert empty
# string so the save step does not turn the valueless <c t="s"/> cells
# into a visible "Name" string (shared string index 0).
$wsProductdata.Range("H2").Value = ""
$wsProductdata.Range("H3").Value = ""
$wsProductdata.Range("H4").Value = ""
$wsProductdata.Range("H5").Value = ""
$wsProductdata.Range("H6").Value = ""
$wsProductdata.Range("H7").Value = ""
$wsProductdata.Range("H8").Value = ""
$wsProductdata.Range("H9").Value = ""
$wsProductdata.Range("H10").Value = ""
$wsProductdata.Range("H11").Value = ""

# --- ForecastedAverageDemand sheet: set demands for last 3 periods (rows 9-11) ---
$wsAvgDemand = $wb.Worksheets.Item("ForecastedAverageDemand")
$wsAvgDemand.Range("B9").Value = 70
$wsAvgDemand.Range("C9").Value = 30
$wsAvgDemand.Range("D9").Value = 50
$wsAvgDemand.Range("E9").Value = 100

$wsAvgDemand.Range("B10").Value = 70
$wsAvgDemand.Range("C10").Value = 30
$wsAvgDemand.Range("D10").Value = 50
$wsAvgDemand.Range("E10").Value = 100

$wsAvgDemand.Range("B11").Value = 70
$wsAvgDemand.Range("C11").Value = 30
$wsAvgDemand.Range("D11").Value = 50
$wsAvgDemand.Range("E11").Value = 100

# --- ForcastedStandardDeviation sheet: set std deviations for last 3 periods (rows 9-11) ---
$wsStdDev = $wb.Worksheets.Item("ForcastedStandardDeviation")
$wsStdDev.Range("B9").Value = 7.166424999999998
$wsStdDev.Range("C9").Value = 3.071324999999999
$wsStdDev.Range("D9").Value = 5.118874999999999
$wsStdDev.Range("E9").Value = 10.23775

$wsStdDev.Range("B10").Value = 8.1997825
$wsStdDev.Range("C10").Value = 3.5141925
$wsStdDev.Range("D10").Value = 5.856987499999999
$wsStdDev.Range("E10").Value = 11.713975

$wsStdDev.Range("B11").Value = 9.129804249999998
$wsStdDev.Range("C11").Value = 3.912773249999999
$wsStdDev.Range("D11").Value = 6.521288749999998
$wsStdDev.Range("E11").Value = 13.0425775
